$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 30: fill in the Read_Mark / Reading grade that were previously blank ---
$ws.Range("H30").Value = 32
$ws.Range("I30").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'

# --- Row 32: newly practiced test (Cambridge 9 Test 4) ---
# Pull row 31's "filled entry" look (thick bottom border, row height) down onto
# row 32 now that it holds real data too.
$ws.Range("C31:L31").Copy()
$ws.Range("C32:L32").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Rows("32").RowHeight = $ws.Rows("31").RowHeight

$ws.Range("D32").Value = 45508
$ws.Range("F32").Value = 33
$ws.Range("G32").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Lis_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'
$ws.Range("H32").Value = 29
$ws.Range("I32").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'
$ws.Range("J32").Value = 1.1000000000000001

# --- Update the active selection to reflect where the user ended up ---
$ws.Range("I39").Select()
